$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.513.24'
$ws.Range("E2").Value = '  -6.21%  '
$ws.Range("D3").Value = '3.279.66'
$ws.Range("E3").Value = '  -6.62%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '550.67'
$ws.Range("E5").Value = '  -4.94%  '
$ws.Range("D6").Value = '181.52'
$ws.Range("E6").Value = '  -5.94%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -4.48%  '
$ws.Range("D9").Value = '3.271.00'
$ws.Range("E9").Value = '  -6.58%  '
$ws.Range("E10").Value = '  -11.30%  '
$ws.Range("E11").Value = '  -6.72%  '
$ws.Range("D12").Value = '46.94'
$ws.Range("E12").Value = '  -8.84%  '
$ws.Range("E13").Value = '  -7.92%  '
$ws.Range("D14").Value = '639.03'
$ws.Range("E14").Value = '  -1.89%  '
$ws.Range("D15").Value = '8.59'
$ws.Range("E15").Value = '  -6.39%  '
$ws.Range("D16").Value = '3.806.08'
$ws.Range("E16").Value = '  -6.19%  '
$ws.Range("E17").Value = '  -2.22%  '
$ws.Range("D18").Value = '65.463.26'
$ws.Range("E18").Value = '  -6.29%  '
$ws.Range("E19").Value = '  -3.18%  '
$ws.Range("D20").Value = '3.280.01'
$ws.Range("E20").Value = '  -6.66%  '
$ws.Range("D21").Value = '11.33'
$ws.Range("E21").Value = '  -8.93%  '
$ws.Range("D22").Value = '0.900'
$ws.Range("E22").Value = '  -5.54%  '
$ws.Range("D23").Value = '17.97'
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("D24").Value = '107.04'
$ws.Range("E24").Value = '  +7.38%  '
$ws.Range("D25").Value = '4.86'
$ws.Range("E25").Value = '  -9.45%  '
$ws.Range("D26").Value = '3.93'
$ws.Range("E26").Value = '  -8.67%  '
$ws.Range("D27").Value = '2.66'
$ws.Range("E27").Value = '  -7.99%  '
$ws.Range("D29").Value = '8.62'
$ws.Range("E29").Value = '  -8.59%  '
$ws.Range("D30").Value = '30.04'
$ws.Range("E30").Value = '  -8.36%  '
$ws.Range("D31").Value = '3.85'
$ws.Range("E31").Value = '  -10.24%  '
$ws.Range("E32").Value = '  -8.22%  '
$ws.Range("E34").Value = '  -5.67%  '
$ws.Range("D35").Value = '3.750.46'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '57.14'
$ws.Range("E37").Value = '  -7.06%  '
$ws.Range("D38").Value = '522.63'
$ws.Range("E38").Value = '  -9.79%  '
$ws.Range("E39").Value = '  -8.98%  '
$ws.Range("D40").Value = '3.32'
$ws.Range("E40").Value = '  -8.17%  '
$ws.Range("E41").Value = '  -4.25%  '
$ws.Range("E42").Value = '  -7.09%  '
$ws.Range("B43").Value = 'CoreDAO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D43").Value = '3.34'
$ws.Range("E43").Value = '  -10.95%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '32.75'
$ws.Range("E44").Value = '  -4.49%  '
$ws.Range("E45").Value = '  -11.21%  '
$ws.Range("D46").Value = '3.27'
$ws.Range("E46").Value = '  -2.46%  '
$ws.Range("E47").Value = '  -7.58%  '
$ws.Range("E48").Value = '  -5.25%  '
$ws.Range("D49").Value = '2.59'
$ws.Range("E49").Value = '  -9.82%  '
$ws.Range("E51").Value = '  +0.95%  '
